$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (4th column) values are being replaced with a "halved" model:
# new value = ceil(old value / 2). Do this in-place for rows 1..109.
for ($r = 1; $r -le 109; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $old = $cell.Value2
    $new = [Math]::Ceiling($old / 2)
    $cell.Value = $new
}

# The new selection left on the sheet is the whole of column E.
$ws.Range("E1:E1048576").Select()
